$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 542.4  # H33: 114.875 -> 542.4
$ws.Cells.Item(33, 9).Value = 499.14285  # I33: 62.714287 -> 499.14285
$ws.Cells.Item(33, 10).Value = 643.3333  # J33: 480 -> 643.3333
$ws.Cells.Item(33, 11).Value = 499.14285  # K33: 62.714287 -> 499.14285
$ws.Cells.Item(33, 12).Value = 643.3333  # L33: 480 -> 643.3333
$ws.Cells.Item(33, 13).Value = -270.14285  # M33: 166.285713 -> -270.14285
$ws.Cells.Item(33, 14).Value = -1101.3333  # N33: -938 -> -1101.3333
$ws.Cells.Item(38, 8).Value = 731  # H38: 374.66666 -> 731
$ws.Cells.Item(38, 10).Value = 1480  # J38: 1000 -> 1480
$ws.Cells.Item(38, 12).Value = 4440  # L38: 3000 -> 4440
$ws.Cells.Item(38, 14).Value = -5184  # N38: -3744 -> -5184
$ws.Cells.Item(40, 8).Value = 2159.3125  # H40: 2575 -> 2159.3125
$ws.Cells.Item(40, 9).Value = 2061  # I40: 2700 -> 2061
$ws.Cells.Item(40, 10).Value = 2285.7144  # J40: 2450 -> 2285.7144
$ws.Cells.Item(40, 11).Value = 2061  # K40: 2700 -> 2061
$ws.Cells.Item(40, 12).Value = 2285.7144  # L40: 2450 -> 2285.7144
$ws.Cells.Item(40, 13).Value = -1886  # M40: -2525 -> -1886
$ws.Cells.Item(40, 14).Value = -2635.7144  # N40: -2800 -> -2635.7144
$ws.Cells.Item(42, 8).Value = 91.333336  # H42: 122.375 -> 91.333336
$ws.Cells.Item(42, 9).Value = 35.75  # I42: 40 -> 35.75
$ws.Cells.Item(42, 10).Value = 135.8  # J42: 171.8 -> 135.8
$ws.Cells.Item(42, 11).Value = 107.25  # K42: 120 -> 107.25
$ws.Cells.Item(42, 12).Value = 407.4  # L42: 515.4000000000001 -> 407.4
$ws.Cells.Item(42, 13).Value = 122.75  # M42: 110 -> 122.75
$ws.Cells.Item(42, 14).Value = -867.4000000000001  # N42: -975.4000000000001 -> -867.4000000000001
$ws.Cells.Item(43, 8).Value = 803.5  # H43: 874.6 -> 803.5
$ws.Cells.Item(43, 9).Value = 634  # I43: 800.5 -> 634
$ws.Cells.Item(43, 10).Value = 849.7273  # J43: 893.125 -> 849.7273
$ws.Cells.Item(43, 11).Value = 634  # K43: 800.5 -> 634
$ws.Cells.Item(43, 12).Value = 849.7273  # L43: 893.125 -> 849.7273
$ws.Cells.Item(43, 13).Value = -565  # M43: -731.5 -> -565
$ws.Cells.Item(43, 14).Value = -987.7273  # N43: -1031.125 -> -987.7273
$ws.Cells.Item(69, 8).Value = 4446625  # H69: 3970411.5 -> 4446625
$ws.Cells.Item(69, 9).Value = 2001.3  # I69: 1993.3077 -> 2001.3
$ws.Cells.Item(69, 11).Value = 6003.9  # K69: 5979.9231 -> 6003.9
$ws.Cells.Item(69, 13).Value = -5129.9  # M69: -5105.9231 -> -5129.9
$ws.Cells.Item(72, 8).Value = 4446625  # H72: 3970411.5 -> 4446625
$ws.Cells.Item(72, 9).Value = 2001.3  # I72: 1993.3077 -> 2001.3
$ws.Cells.Item(72, 11).Value = 18011.7  # K72: 17939.7693 -> 18011.7
$ws.Cells.Item(72, 13).Value = -13643.7  # M72: -13571.7693 -> -13643.7
$ws.Cells.Item(96, 8).Value = 439.85715  # H96: 491.66666 -> 439.85715
$ws.Cells.Item(96, 10).Value = 326.33334  # J96: 425 -> 326.33334
$ws.Cells.Item(96, 12).Value = 979.0000200000001  # L96: 1275 -> 979.0000200000001
$ws.Cells.Item(96, 14).Value = -3725.00002  # N96: -4021 -> -3725.00002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 5000250  # H5: 2000200.6 -> 5000250
$ws.Cells.Item(5, 9).Value = 5000250  # I5: 3333540.2 -> 5000250
$ws.Cells.Item(5, 10).Value = 0  # J5: 191 -> 0
$ws.Cells.Item(5, 11).Value = 5000250  # K5: 3333540.2 -> 5000250
$ws.Cells.Item(5, 12).Value = 0  # L5: 191 -> 0
$ws.Cells.Item(5, 13).Value = -5000138  # M5: -3333428.2 -> -5000138
$ws.Cells.Item(5, 14).ClearContents()  # N5: -415 -> (cleared)
$ws.Cells.Item(74, 8).Value = 9520.375  # H74: 10022.4 -> 9520.375
$ws.Cells.Item(74, 9).Value = 2275.0908  # I74: 2276 -> 2275.0908
$ws.Cells.Item(74, 10).Value = 25460  # J74: 31325 -> 25460
$ws.Cells.Item(74, 11).Value = 2275.0908  # K74: 2276 -> 2275.0908
$ws.Cells.Item(74, 12).Value = 25460  # L74: 31325 -> 25460
$ws.Cells.Item(74, 13).Value = -1401.0908  # M74: -1402 -> -1401.0908
$ws.Cells.Item(74, 14).Value = -27208  # N74: -33073 -> -27208
$ws.Cells.Item(77, 8).Value = 9520.375  # H77: 10022.4 -> 9520.375
$ws.Cells.Item(77, 9).Value = 2275.0908  # I77: 2276 -> 2275.0908
$ws.Cells.Item(77, 10).Value = 25460  # J77: 31325 -> 25460
$ws.Cells.Item(77, 11).Value = 11375.454  # K77: 11380 -> 11375.454
$ws.Cells.Item(77, 12).Value = 127300  # L77: 156625 -> 127300
$ws.Cells.Item(77, 13).Value = -7007.454  # M77: -7012 -> -7007.454
$ws.Cells.Item(77, 14).Value = -136036  # N77: -165361 -> -136036
$ws.Cells.Item(97, 8).Value = 33343850  # H97: 18524556 -> 33343850
$ws.Cells.Item(97, 9).Value = 41679624  # I97: 19614176 -> 41679624
$ws.Cells.Item(97, 10).Value = 755  # J97: 1000 -> 755
$ws.Cells.Item(97, 11).Value = 41679624  # K97: 19614176 -> 41679624
$ws.Cells.Item(97, 12).Value = 755  # L97: 1000 -> 755
$ws.Cells.Item(97, 13).Value = -41679128  # M97: -19613680 -> -41679128
$ws.Cells.Item(97, 14).Value = -1747  # N97: -1992 -> -1747
$ws.Cells.Item(102, 8).Value = 976.6667  # H102: 2237.3333 -> 976.6667
$ws.Cells.Item(102, 9).Value = 976.6667  # I102: 2390.625 -> 976.6667
$ws.Cells.Item(102, 10).Value = 0  # J102: 1011 -> 0
$ws.Cells.Item(102, 11).Value = 976.6667  # K102: 2390.625 -> 976.6667
$ws.Cells.Item(102, 12).Value = 0  # L102: 1011 -> 0
$ws.Cells.Item(102, 13).Value = 645.3333  # M102: -768.625 -> 645.3333
$ws.Cells.Item(102, 14).ClearContents()  # N102: -4255 -> (cleared)
$ws.Cells.Item(122, 8).Value = 2268.3845  # H122: 2657.7 -> 2268.3845
$ws.Cells.Item(122, 9).Value = 2548.9  # I122: 2947.125 -> 2548.9
$ws.Cells.Item(122, 10).Value = 1333.3334  # J122: 1500 -> 1333.3334
$ws.Cells.Item(122, 11).Value = 7646.700000000001  # K122: 8841.375 -> 7646.700000000001
$ws.Cells.Item(122, 12).Value = 4000.0002  # L122: 4500 -> 4000.0002
$ws.Cells.Item(122, 13).Value = -5196.700000000001  # M122: -6391.375 -> -5196.700000000001
$ws.Cells.Item(122, 14).Value = -8900.0002  # N122: -9400 -> -8900.0002
$ws.Cells.Item(132, 8).Value = 2827.48  # H132: 2820.25 -> 2827.48
$ws.Cells.Item(132, 9).Value = 2202.4736  # I132: 2255.0527 -> 2202.4736
$ws.Cells.Item(132, 10).Value = 4806.6665  # J132: 4968 -> 4806.6665
$ws.Cells.Item(132, 11).Value = 6607.4208  # K132: 6765.158100000001 -> 6607.4208
$ws.Cells.Item(132, 12).Value = 14419.9995  # L132: 14904 -> 14419.9995
$ws.Cells.Item(132, 13).Value = -4077.4208  # M132: -4235.158100000001 -> -4077.4208
$ws.Cells.Item(132, 14).Value = -19479.9995  # N132: -19964 -> -19479.9995
$ws.Cells.Item(139, 8).Value = 0  # H139: 60715 -> 0
$ws.Cells.Item(139, 10).Value = 0  # J139: 60715 -> 0
$ws.Cells.Item(139, 12).Value = 0  # L139: 60715 -> 0
$ws.Cells.Item(139, 14).ClearContents()  # N139: -70995 -> (cleared)

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 5000250  # H4: 2000200.6 -> 5000250
$ws.Cells.Item(4, 9).Value = 5000250  # I4: 3333540.2 -> 5000250
$ws.Cells.Item(4, 10).Value = 0  # J4: 191 -> 0
$ws.Cells.Item(4, 11).Value = 5000250  # K4: 3333540.2 -> 5000250
$ws.Cells.Item(4, 12).Value = 0  # L4: 191 -> 0
$ws.Cells.Item(4, 13).Value = -5000135  # M4: -3333425.2 -> -5000135
$ws.Cells.Item(4, 14).ClearContents()  # N4: -421 -> (cleared)
$ws.Cells.Item(29, 8).Value = 2508  # H29: 2500 -> 2508
$ws.Cells.Item(29, 9).Value = 2508  # I29: 2500 -> 2508
$ws.Cells.Item(29, 11).Value = 2508  # K29: 2500 -> 2508
$ws.Cells.Item(29, 13).Value = -2219  # M29: -2211 -> -2219
$ws.Cells.Item(36, 8).Value = 1290  # H36: 1426.25 -> 1290
$ws.Cells.Item(36, 9).Value = 701.25  # I36: 772.8570999999999 -> 701.25
$ws.Cells.Item(36, 11).Value = 701.25  # K36: 772.8570999999999 -> 701.25
$ws.Cells.Item(36, 13).Value = -167.25  # M36: -238.8570999999999 -> -167.25
$ws.Cells.Item(75, 8).Value = 83473.09  # H75: 75289.46000000001 -> 83473.09
$ws.Cells.Item(75, 9).Value = 4746  # I75: 5391.8887 -> 4746
$ws.Cells.Item(75, 10).Value = 293412  # J75: 232559 -> 293412
$ws.Cells.Item(75, 11).Value = 4746  # K75: 5391.8887 -> 4746
$ws.Cells.Item(75, 12).Value = 293412  # L75: 232559 -> 293412
$ws.Cells.Item(75, 13).Value = -3810  # M75: -4455.8887 -> -3810
$ws.Cells.Item(75, 14).Value = -295284  # N75: -234431 -> -295284
$ws.Cells.Item(78, 8).Value = 83473.09  # H78: 75289.46000000001 -> 83473.09
$ws.Cells.Item(78, 9).Value = 4746  # I78: 5391.8887 -> 4746
$ws.Cells.Item(78, 10).Value = 293412  # J78: 232559 -> 293412
$ws.Cells.Item(78, 11).Value = 14238  # K78: 16175.6661 -> 14238
$ws.Cells.Item(78, 12).Value = 880236  # L78: 697677 -> 880236
$ws.Cells.Item(78, 13).Value = -9558  # M78: -11495.6661 -> -9558
$ws.Cells.Item(78, 14).Value = -889596  # N78: -707037 -> -889596
$ws.Cells.Item(99, 8).Value = 1390  # H99: 1481.75 -> 1390
$ws.Cells.Item(99, 9).Value = 1513.75  # I99: 1481.75 -> 1513.75
$ws.Cells.Item(99, 10).Value = 400  # J99: 0 -> 400
$ws.Cells.Item(99, 11).Value = 1513.75  # K99: 1481.75 -> 1513.75
$ws.Cells.Item(99, 12).Value = 400  # L99: 0 -> 400
$ws.Cells.Item(99, 13).Value = -15.75  # M99: 16.25 -> -15.75
$ws.Cells.Item(99, 14).Value = -3396  # N99: None -> -3396
$ws.Cells.Item(132, 8).Value = 12498.333  # H132: 12661.667 -> 12498.333
$ws.Cells.Item(132, 10).Value = 12498.333  # J132: 12661.667 -> 12498.333
$ws.Cells.Item(132, 12).Value = 12498.333  # L132: 12661.667 -> 12498.333
$ws.Cells.Item(132, 14).Value = -22618.333  # N132: -22781.667 -> -22618.333
$ws.Cells.Item(140, 8).Value = 50500  # H140: 50780 -> 50500
$ws.Cells.Item(140, 10).Value = 50500  # J140: 50780 -> 50500
$ws.Cells.Item(140, 12).Value = 50500  # L140: 50780 -> 50500
$ws.Cells.Item(140, 14).Value = -60860  # N140: -61140 -> -60860

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 0  # H6: 24824.4 -> 0
$ws.Cells.Item(6, 9).Value = 0  # I6: 40704 -> 0
$ws.Cells.Item(6, 10).Value = 0  # J6: 1005 -> 0
$ws.Cells.Item(6, 11).Value = 0  # K6: 40704 -> 0
$ws.Cells.Item(6, 12).Value = 0  # L6: 1005 -> 0
$ws.Cells.Item(6, 13).ClearContents()  # M6: -40591 -> (cleared)
$ws.Cells.Item(6, 14).ClearContents()  # N6: -1231 -> (cleared)
$ws.Cells.Item(7, 8).Value = 54.307693  # H7: 98.625 -> 54.307693
$ws.Cells.Item(7, 9).Value = 35.42857  # I7: 69.8 -> 35.42857
$ws.Cells.Item(7, 10).Value = 76.333336  # J7: 146.66667 -> 76.333336
$ws.Cells.Item(7, 11).Value = 35.42857  # K7: 69.8 -> 35.42857
$ws.Cells.Item(7, 12).Value = 76.333336  # L7: 146.66667 -> 76.333336
$ws.Cells.Item(7, 13).Value = 77.57142999999999  # M7: 43.2 -> 77.57142999999999
$ws.Cells.Item(7, 14).Value = -302.333336  # N7: -372.66667 -> -302.333336
$ws.Cells.Item(17, 8).Value = 14004.75  # H17: 0 -> 14004.75
$ws.Cells.Item(17, 9).Value = 17670  # I17: 0 -> 17670
$ws.Cells.Item(17, 10).Value = 3009  # J17: 0 -> 3009
$ws.Cells.Item(17, 11).Value = 17670  # K17: 0 -> 17670
$ws.Cells.Item(17, 12).Value = 3009  # L17: 0 -> 3009
$ws.Cells.Item(17, 13).Value = -17496  # M17: None -> -17496
$ws.Cells.Item(17, 14).Value = -3357  # N17: None -> -3357
$ws.Cells.Item(25, 8).Value = 6500  # H25: 2400 -> 6500
$ws.Cells.Item(25, 9).Value = 0  # I25: 300 -> 0
$ws.Cells.Item(25, 10).Value = 6500  # J25: 3100 -> 6500
$ws.Cells.Item(25, 11).Value = 0  # K25: 300 -> 0
$ws.Cells.Item(25, 12).Value = 6500  # L25: 3100 -> 6500
$ws.Cells.Item(25, 13).ClearContents()  # M25: -126 -> (cleared)
$ws.Cells.Item(25, 14).Value = -6848  # N25: -3448 -> -6848
$ws.Cells.Item(51, 8).Value = 27933.334  # H51: 29400 -> 27933.334
$ws.Cells.Item(51, 10).Value = 27933.334  # J51: 29400 -> 27933.334
$ws.Cells.Item(51, 12).Value = 27933.334  # L51: 29400 -> 27933.334
$ws.Cells.Item(51, 14).Value = -29405.334  # N51: -30872 -> -29405.334
$ws.Cells.Item(61, 8).Value = 27933.334  # H61: 29400 -> 27933.334
$ws.Cells.Item(61, 10).Value = 27933.334  # J61: 29400 -> 27933.334
$ws.Cells.Item(61, 12).Value = 27933.334  # L61: 29400 -> 27933.334
$ws.Cells.Item(61, 14).Value = -28629.334  # N61: -30096 -> -28629.334

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(42, 8).Value = 29700  # H42: 14263 -> 29700
$ws.Cells.Item(42, 9).Value = 0  # I42: 14263 -> 0
$ws.Cells.Item(42, 10).Value = 29700  # J42: 0 -> 29700
$ws.Cells.Item(42, 11).Value = 0  # K42: 14263 -> 0
$ws.Cells.Item(42, 12).Value = 29700  # L42: 0 -> 29700
$ws.Cells.Item(42, 13).ClearContents()  # M42: -13778 -> (cleared)
$ws.Cells.Item(42, 14).Value = -30670  # N42: None -> -30670
$ws.Cells.Item(97, 8).Value = 91672.63  # H97: 125750 -> 91672.63
$ws.Cells.Item(97, 9).Value = 125712.375  # I97: 143571.42 -> 125712.375
$ws.Cells.Item(97, 10).Value = 900  # J97: 1000 -> 900
$ws.Cells.Item(97, 11).Value = 125712.375  # K97: 143571.42 -> 125712.375
$ws.Cells.Item(97, 12).Value = 900  # L97: 1000 -> 900
$ws.Cells.Item(97, 13).Value = -125216.375  # M97: -143075.42 -> -125216.375
$ws.Cells.Item(97, 14).Value = -1892  # N97: -1992 -> -1892
$ws.Cells.Item(115, 8).Value = 29700  # H115: 14263 -> 29700
$ws.Cells.Item(115, 9).Value = 0  # I115: 14263 -> 0
$ws.Cells.Item(115, 10).Value = 29700  # J115: 0 -> 29700
$ws.Cells.Item(115, 11).Value = 0  # K115: 14263 -> 0
$ws.Cells.Item(115, 12).Value = 29700  # L115: 0 -> 29700
$ws.Cells.Item(115, 13).ClearContents()  # M115: -13088 -> (cleared)
$ws.Cells.Item(115, 14).Value = -32050  # N115: None -> -32050
$ws.Cells.Item(138, 8).Value = 53500  # H138: 60500 -> 53500
$ws.Cells.Item(138, 10).Value = 53500  # J138: 60500 -> 53500
$ws.Cells.Item(138, 12).Value = 53500  # L138: 60500 -> 53500
$ws.Cells.Item(138, 14).Value = -63780  # N138: -70780 -> -63780

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 4600  # H9: 1399.75 -> 4600
$ws.Cells.Item(9, 9).Value = 400  # I9: 933 -> 400
$ws.Cells.Item(9, 10).Value = 6000  # J9: 2800 -> 6000
$ws.Cells.Item(9, 11).Value = 400  # K9: 933 -> 400
$ws.Cells.Item(9, 12).Value = 6000  # L9: 2800 -> 6000
$ws.Cells.Item(9, 13).Value = -176  # M9: -709 -> -176
$ws.Cells.Item(9, 14).Value = -6448  # N9: -3248 -> -6448
$ws.Cells.Item(55, 8).Value = 570.7143  # H55: 372.54544 -> 570.7143
$ws.Cells.Item(55, 9).Value = 558.6  # I55: 385.14285 -> 558.6
$ws.Cells.Item(55, 10).Value = 601  # J55: 350.5 -> 601
$ws.Cells.Item(55, 11).Value = 558.6  # K55: 385.14285 -> 558.6
$ws.Cells.Item(55, 12).Value = 601  # L55: 350.5 -> 601
$ws.Cells.Item(55, 13).Value = -385.6  # M55: -212.14285 -> -385.6
$ws.Cells.Item(55, 14).Value = -947  # N55: -696.5 -> -947

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4566.6665  # H81: 4700 -> 4566.6665
$ws.Cells.Item(81, 9).Value = 2000  # I81: 1700 -> 2000
$ws.Cells.Item(81, 10).Value = 4961.5386  # J81: 4900 -> 4961.5386
$ws.Cells.Item(81, 11).Value = 4000  # K81: 3400 -> 4000
$ws.Cells.Item(81, 12).Value = 9923.0772  # L81: 9800 -> 9923.0772
$ws.Cells.Item(81, 13).Value = -2939  # M81: -2339 -> -2939
$ws.Cells.Item(81, 14).Value = -12045.0772  # N81: -11922 -> -12045.0772
$ws.Cells.Item(84, 8).Value = 4566.6665  # H84: 4700 -> 4566.6665
$ws.Cells.Item(84, 9).Value = 2000  # I84: 1700 -> 2000
$ws.Cells.Item(84, 10).Value = 4961.5386  # J84: 4900 -> 4961.5386
$ws.Cells.Item(84, 11).Value = 20000  # K84: 17000 -> 20000
$ws.Cells.Item(84, 12).Value = 49615.386  # L84: 49000 -> 49615.386
$ws.Cells.Item(84, 13).Value = -14696  # M84: -11696 -> -14696
$ws.Cells.Item(84, 14).Value = -60223.386  # N84: -59608 -> -60223.386
$ws.Cells.Item(100, 8).Value = 501  # H100: 312.4 -> 501
$ws.Cells.Item(100, 9).Value = 202  # I100: 190.5 -> 202
$ws.Cells.Item(100, 11).Value = 404  # K100: 381 -> 404
$ws.Cells.Item(100, 13).Value = 137  # M100: 160 -> 137
